$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the hidden "_GoBack" bookmark so it sits right after the
#    "Programming Plan" run in the first paragraph (instead of in the
#    middle of the URL run further down).
#
#    Locate "Programming Plan" with Find so we don't depend on a
#    hard-coded offset, then collapse to the spot right after it.
#
#    Word's real behaviour treats a bookmark collapsed exactly at the
#    end of a paragraph's content as spanning the *whole* paragraph
#    (bookmarkStart jumps to the paragraph start, bookmarkEnd spills
#    into the next paragraph), so instead we briefly insert a
#    throwaway character right after "Plan", anchor the collapsed
#    bookmark at that (now mid-paragraph) position, and then remove
#    the throwaway character again.
# ---------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute("Programming Plan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterTitle = $titleRange.End

$tempMarker = $d.Range($afterTitle, $afterTitle)
$tempMarker.Text = "X"

$bookmarkSpot = $d.Range($afterTitle, $afterTitle)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$cleanup = $d.Range($afterTitle, $afterTitle + 1)
$cleanup.Text = ""

# ---------------------------------------------------------------------
# 2) The last paragraph (the Google Drive URL) previously held the
#    bookmark mid-run, splitting the URL text into two separate runs:
#    "https://drive.google.com/drive/folders/0BxEiYrSnB7tgTHN" and
#    "6Vk04Sm9QUkE". Adding the bookmark above already removed it from
#    here; now collapse the two runs back into a single contiguous run
#    with the full URL text.
# ---------------------------------------------------------------------
$urlPara = $d.Paragraphs.Last
$urlStart = $urlPara.Range.Start
$urlEnd = $urlPara.Range.End - 1

$urlRange = $d.Range($urlStart, $urlEnd)
$urlRange.Text = ""

$urlInsert = $d.Range($urlStart, $urlStart)
$urlInsert.Text = "https://drive.google.com/drive/folders/0BxEiYrSnB7tgTHN6Vk04Sm9QUkE"

Write-Output "Applied edits: relocated _GoBack bookmark and merged URL runs."
